$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '55.534.46'
$ws.Range('E2').Value = '  +5.36%  '
$ws.Range('D3').Value = '2.501.54'
$ws.Range('E3').Value = '  +9.97%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '''482.29'
$ws.Range('E5').Value = '  +11.63%  '
$ws.Range('D6').Value = '''140.34'
$ws.Range('E6').Value = '  +17.80%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +8.34%  '
$ws.Range('D9').Value = '2.499.73'
$ws.Range('E9').Value = '  +9.36%  '
$ws.Range('D10').Value = '''0.0986'
$ws.Range('E10').Value = '  +8.98%  '
$ws.Range('D11').Value = '''5.47'
$ws.Range('E11').Value = '  +4.00%  '
$ws.Range('E12').Value = '  +7.37%  '
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').Value = '2.933.82'
$ws.Range('E14').Value = '  +8.48%  '
$ws.Range('D15').Value = '55.546.71'
$ws.Range('E15').Value = '  +5.40%  '
$ws.Range('D16').Value = '''20.66'
$ws.Range('E16').Value = '  +10.34%  '
$ws.Range('E17').Value = '  +15.93%  '
$ws.Range('D18').Value = '2.496.50'
$ws.Range('E18').Value = '  +7.49%  '
$ws.Range('E19').Value = '  +11.19%  '
$ws.Range('D20').Value = '''320.28'
$ws.Range('E20').Value = '  +8.05%  '
$ws.Range('D21').Value = '''10.04'
$ws.Range('E21').Value = '  +10.23%  '
$ws.Range('D22').Value = '''0.998'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '''5.70'
$ws.Range('E23').Value = '  +7.09%  '
$ws.Range('D24').Value = '''57.82'
$ws.Range('E24').Value = '  +5.38%  '
$ws.Range('D25').Value = '''0.167'
$ws.Range('E25').Value = '  +11.66%  '
$ws.Range('D26').Value = '''0.409'
$ws.Range('E26').Value = '  +12.67%  '
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').Value = '2.604.65'
$ws.Range('E28').Value = '  +8.35%  '
$ws.Range('D29').Value = '''7.37'
$ws.Range('E29').Value = '  +6.81%  '
$ws.Range('D30').Value = '0.0₃0797'
$ws.Range('E30').Value = '  +16.71%  '
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').Value = '''149.40'
$ws.Range('E32').Value = '  +3.82%  '
$ws.Range('D33').Value = '''18.17'
$ws.Range('E33').Value = '  +7.19%  '
$ws.Range('E34').Value = '  +11.99%  '
$ws.Range('D35').Value = '''5.16'
$ws.Range('E35').Value = '  +11.40%  '
$ws.Range('D36').Value = '''3.72'
$ws.Range('E36').Value = '  +6.25%  '
$ws.Range('E37').Value = '  +12.78%  '
$ws.Range('D38').Value = '''0.861'
$ws.Range('E38').Value = '  +6.85%  '
$ws.Range('D39').Value = '''34.19'
$ws.Range('E39').Value = '  +3.43%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '''0.608'
$ws.Range('E41').Value = '  +18.89%  '
$ws.Range('D42').Value = '''0.0554'
$ws.Range('E42').Value = '  +12.96%  '
$ws.Range('D43').Value = '''3.41'
$ws.Range('E43').Value = '  +9.58%  '
$ws.Range('E44').Value = '  +10.74%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.979.79'
$ws.Range('E45').Value = '  +4.07%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '''10.13'
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('D47').Value = '''0.0903'
$ws.Range('E47').Value = '  +10.29%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''4.62'
$ws.Range('E48').Value = '  +20.09%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '''0.0223'
$ws.Range('E49').Value = '  +8.50%  '
$ws.Range('D50').Value = '''251.37'
$ws.Range('E50').Value = '  +35.78%  '
$ws.Range('D51').Value = '''17.60'
$ws.Range('E51').Value = '  +12.47%  '